$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation ("Cristal" variety, Región del Maule) was
# recorded for this market. It belongs chronologically right after the
# existing row 14, so insert a fresh row at position 15 - this pushes the
# former rows 15..84 down to 16..85 - and then fill in the new record.
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(15, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15, 3).Value = "Ñuble"
$ws.Cells.Item(15, 4).Value = 44687
$ws.Cells.Item(15, 5).Value = 16
$ws.Cells.Item(15, 6).Value = 100112021
$ws.Cells.Item(15, 7).Value = "Ají"
$ws.Cells.Item(15, 8).Value = "Cristal"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 24000
$ws.Cells.Item(15, 12).Value = 25000
$ws.Cells.Item(15, 13).Value = 24500
$ws.Cells.Item(15, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región del Maule"
$ws.Cells.Item(15, 16).Value = 980
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"
